$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 1256-1257, pushing the existing data (previously
# rows 1256:1315) down to 1258:1317. Excel's Insert copies the formatting
# (e.g. the date style on column D) from the row above automatically.
$ws.Rows.Item(1256).Resize(2).Insert()

# New weekly record #1 (row 1256) - Coliflor, Primera, 2023-08-09
$ws.Range("A1256").Value = 6
$ws.Range("B1256").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1256").Value = "Metropolitana"
$ws.Range("D1256").Value = 45147
$ws.Range("E1256").Value = 13
$ws.Range("F1256").Value = 100112008
$ws.Range("G1256").Value = "Coliflor"
$ws.Range("H1256").Value = "Sin especificar"
$ws.Range("I1256").Value = "Primera"
$ws.Range("J1256").Value = 6300
$ws.Range("K1256").Value = 600
$ws.Range("L1256").Value = 700
$ws.Range("M1256").Value = 660
$ws.Range("N1256").Value = "`$/unidad"
$ws.Range("O1256").Value = "Región Metropolitana"
$ws.Range("P1256").Value = 660
$ws.Range("Q1256").Value = 1
$ws.Range("R1256").Value = "Hortaliza"

# New weekly record #2 (row 1257) - Coliflor, Segunda, 2023-08-09
$ws.Range("A1257").Value = 6
$ws.Range("B1257").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1257").Value = "Metropolitana"
$ws.Range("D1257").Value = 45147
$ws.Range("E1257").Value = 13
$ws.Range("F1257").Value = 100112008
$ws.Range("G1257").Value = "Coliflor"
$ws.Range("H1257").Value = "Sin especificar"
$ws.Range("I1257").Value = "Segunda"
$ws.Range("J1257").Value = 4500
$ws.Range("K1257").Value = 400
$ws.Range("L1257").Value = 500
$ws.Range("M1257").Value = 433
$ws.Range("N1257").Value = "`$/unidad"
$ws.Range("O1257").Value = "Región Metropolitana"
$ws.Range("P1257").Value = 433
$ws.Range("Q1257").Value = 1
$ws.Range("R1257").Value = "Hortaliza"

Write-Host "Inserted 2 new rows; dimension now" $ws.UsedRange.Rows.Count
